$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price/Volume columns so values that look numeric
# (e.g. "1.000", "305.97") are stored as literal text, matching the source data
# (which uses inline strings, not numbers, for these columns).
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range('D2').Value = '23.463.33'
$ws.Range('E2').Value = '  +1.18%  '
$ws.Range('D3').Value = '1.638.17'
$ws.Range('E3').Value = '  +2.25%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').Value = '305.97'
$ws.Range('E6').Value = '  +0.86%  '
$ws.Range('D7').Value = '0.3757'
$ws.Range('E7').Value = '  -0.68%  '
$ws.Range('D8').Value = '52.01'
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  +0.51%  '
$ws.Range('E10').Value = '  -0.91%  '
$ws.Range('D11').Value = '0.08143'
$ws.Range('D13').Value = '22.95'
$ws.Range('E13').Value = '  +0.75%  '
$ws.Range('D14').Value = '6.631'
$ws.Range('E14').Value = '  +0.44%  '
$ws.Range('E15').Value = '  +2.45%  '
$ws.Range('D16').Value = '7.367'
$ws.Range('E16').Value = '  -0.74%  '
$ws.Range('D17').Value = '1.632.43'
$ws.Range('E17').Value = '  +1.83%  '
$ws.Range('D18').Value = '94.52'
$ws.Range('E18').Value = '  +0.56%  '
$ws.Range('D19').Value = '0.06907'
$ws.Range('E19').Value = '  +0.43%  '
$ws.Range('E20').Value = '  +0.53%  '
$ws.Range('D21').Value = '6.537'
$ws.Range('E21').Value = '  -0.21%  '
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').Value = '23.460.25'
$ws.Range('E23').Value = '  +1.21%  '
$ws.Range('D24').Value = '12.77'
$ws.Range('E24').Value = '  -1.61%  '
$ws.Range('D25').Value = '3.069'
$ws.Range('E25').Value = '  +3.17%  '
$ws.Range('E26').Value = '  +0.92%  '
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('D28').Value = '151.18'
$ws.Range('E28').Value = '  +1.23%  '
$ws.Range('D29').Value = '5.328'
$ws.Range('E29').Value = '  +1.46%  '
$ws.Range('D30').Value = '136.78'
$ws.Range('E30').Value = '  +2.13%  '
$ws.Range('D31').Value = '2.309'
$ws.Range('E31').Value = '  -2.19%  '
$ws.Range('D32').Value = '1.815.45'
$ws.Range('E32').Value = '  +1.91%  '
$ws.Range('D33').Value = '6.755'
$ws.Range('E33').Value = '  -0.36%  '
$ws.Range('D34').Value = '0.9612'
$ws.Range('E34').Value = '  -0.99%  '
$ws.Range('D35').Value = '0.02839'
$ws.Range('E35').Value = '  +4.39%  '
$ws.Range('E36').Value = '  +0.42%  '
$ws.Range('D37').Value = '0.07309'
$ws.Range('E37').Value = '  -2.81%  '
$ws.Range('D38').Value = '0.2528'
$ws.Range('E38').Value = '  +0.95%  '
$ws.Range('D39').Value = '0.08840'
$ws.Range('E39').Value = '  +0.42%  '
$ws.Range('D40').Value = '6.123'
$ws.Range('E40').Value = '  +0.77%  '
$ws.Range('D41').Value = '1.376'
$ws.Range('E41').Value = '  +1.06%  '
$ws.Range('D42').Value = '0.7094'
$ws.Range('E42').Value = '  -0.27%  '
$ws.Range('D43').Value = '12.47'
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('D44').Value = '16.11'
$ws.Range('E44').Value = '  +3.30%  '
$ws.Range('D45').Value = '0.6552'
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('E46').Value = '  +1.00%  '
$ws.Range('D47').Value = '1.000'
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('D48').Value = '4.010'
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('D49').Value = '0.07972'
$ws.Range('E49').Value = '  +0.24%  '
$ws.Range('D50').Value = '128.65'
$ws.Range('E50').Value = '  -2.77%  '
$ws.Range('D51').Value = '1.204'
$ws.Range('E51').Value = '  +0.16%  '

# Restore the original (default/no explicit) cell formatting so the only
# change versus the source file is the text content itself.
$rng.ClearFormats()
